{"js": "// Remove the \"Word version of this document\" list item (a hyperlink to the\n// .docx version) from the \"Additional resources\" list \u2014 a PDF version is\n// being added elsewhere, so the in-document link to the raw Word file goes\n// away.\nconst results = context.document.body.search(\"Word version of this document\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  const para = results.items[i].paragraphs.getFirst();\n  para.delete();\n}\nawait context.sync();\n", "ps1": "# Remove the \"Word version of this document\" list item (a hyperlink to the\n# .docx version) from the \"Additional resources\" list \u2014 a PDF version is\n# being added elsewhere, so the in-document link to the raw Word file goes\n# away.\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$found = $range.Find.Execute(\"Word version of this document\")\n\nif ($found) {\n    $para = $range.Paragraphs(1)\n    $para.Range.Delete()\n}\n"}
